{"js": "// Insert a new paragraph \"Hello2\" right after the paragraph containing\n// \"Hello\", before the existing (empty) paragraph that follows it.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Find the paragraph whose text is \"Hello\" (the first paragraph in the doc).\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === \"Hello\") {\n    target = paragraphs.items[i];\n    break;\n  }\n}\nif (!target) {\n  target = paragraphs.items[0];\n}\n\n// Insert the new paragraph with text \"Hello2\" directly after it.\ntarget.insertParagraph(\"Hello2\", Word.InsertLocation.after);\nawait context.sync();\n", "ps1": "# Insert a new paragraph \"Hello2\" right after the paragraph containing\n# \"Hello\", before the existing (empty) paragraph that follows it.\n$d = $word.ActiveDocument\n\n# Locate the paragraph whose text is \"Hello\" (strip the trailing\n# paragraph-mark character that Range.Text always carries).\n$target = $null\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.Text.TrimEnd(\"`r\") -eq \"Hello\") {\n        $target = $p\n        break\n    }\n}\nif ($target -eq $null) {\n    $target = $d.Paragraphs.Item(1)\n}\n\n# Insert a brand-new paragraph mark right after it, then fill it with text.\n$target.Range.InsertParagraphAfter()\n$newPara = $target.Next()\n$newPara.Range.InsertBefore(\"Hello2\")\n"}
